# Textgrids Speak pt 9, Praat file
# Adds four new participants (rows 12-14 "filled in", rows 15-16 extended) to the
# "Participants" sheet, expands the note in I10, and moves the active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Participants")
$ws.Activate()

# --- Expand the incidence note for participant 9 (row 10) ---
$ws.Range("I10").Value = "Speaks fast. Tends to ¨cut¨ syllables in the sense that they are not heard clearly at all times. Typed-said the colour instead of the word in the first couple of trials in Practice 1. I corrected him and he did the rest correctly."

# --- Rows 12-14: already have Code/Version/List/Language_test values but are missing
# --- Name/ID/Date and the shaded formatting that the rest of the filled rows use.
# --- Copy the formatting from the fully-filled row 11 so styles match exactly.
$ws.Range("A11:H11").Copy()
$ws.Range("A12:H12").PasteSpecial(-4122)
$ws.Range("A13:H13").PasteSpecial(-4122)
$ws.Range("A14:H14").PasteSpecial(-4122)

$ws.Range("B12").Value = "Irene Pereira López"
$ws.Range("C12").Value = 11570
$ws.Range("H12").Value = 45257

$ws.Range("B13").Value = "Ana Fernández Rubio"
$ws.Range("C13").Value = 9754
$ws.Range("H13").Value = 45257

$ws.Range("B14").Value = "Nerea Pérez Arriazu"
$ws.Range("C14").Value = 10344
$ws.Range("H14").Value = 45257

# --- Rows 15-16: new Name/ID/Date values, unshaded (same look as the rest of the row). ---
$ws.Range("B15").Value = "Haizea Lavega Torrado"
$ws.Range("C15").Value = 11614
$ws.Range("H15").Value = 45258
$ws.Range("H15").NumberFormat = "mm-dd-yy"

$ws.Range("B16").Value = "Clara Lorenzo Tabueña"
$ws.Range("C16").Value = 11782
$ws.Range("H15").Copy()
$ws.Range("H16").PasteSpecial(-4122)
$ws.Range("H16").Value = 45258

# --- View state: scroll back to the top and move the active selection to H17. ---
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
$ws.Range("H17").Select()
